$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @{
    2  = @(0.3731803534348715, 1.522700897529608, 4.158686648532057, 2.039285818254042, 2.032504295995627, 37)
    3  = @(-0.1240574867935324, 1.460950150426982, 3.796378227175162, 1.948429682378905, 1.972058872036625, 36)
    4  = @(0.4934895688106978, 1.327481032087486, 3.305097368472943, 1.817992675582865, 1.77527785305137, 35)
    5  = @(-0.08377265579416678, 1.555732842827626, 4.194765129124354, 2.048112577258475, 2.077173248366187, 34)
    6  = @(0.4363569993891649, 1.325938237731197, 3.276676303256219, 1.810159192793888, 1.784016472626006, 33)
    7  = @(0.09168004481066386, 1.405118573146903, 3.289626680765243, 1.813732803023985, 1.840398680552495, 32)
    8  = @(0.4864221507894972, 1.346678349091636, 3.457793333901256, 1.859514273648163, 1.824434081926045, 31)
    9  = @(0.00377647522211878, 1.4050466251974, 3.34404499931418, 1.828673015963811, 1.85993068212021, 30)
    10 = @(0.5773184669214476, 1.354169541096897, 3.547916642886769, 1.883591421430552, 1.824671994873412, 29)
    11 = @(-0.06534739717242485, 1.454355996855039, 3.515063405020061, 1.874850235357497, 1.908093943603871, 28)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
